$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Modelo" header column (F1), reusing the same formatting as the
# existing header cells (bold font, centered, bordered) by copying the
# format from the neighboring header cell E1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "Modelo"

# New model description value (F2), no special formatting.
$ws.Range("F2").Value = "Pipeline(steps=[('model', LinearRegression())])"
